# Weekly update: insert two new weekly price records (row 26 and 27)
# for "Primera" and "Segunda" quality Sandia, pushing the previously
# existing rows (old rows 26-62) down by two positions (new rows 28-64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 26, shifting everything below down.
$ws.Rows("26:27").Insert()

# New row 26: Sandia, Primera, date 2021-12-09 (serial 44539)
$ws.Cells.Item(26,1).Value  = 8
$ws.Cells.Item(26,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(26,3).Value  = "Coquimbo"
$ws.Cells.Item(26,4).Value  = 44539
$ws.Cells.Item(26,5).Value  = 4
$ws.Cells.Item(26,6).Value  = 100112028
$ws.Cells.Item(26,7).Value  = "Sandia"
$ws.Cells.Item(26,8).Value  = "Sin especificar"
$ws.Cells.Item(26,9).Value  = "Primera"
$ws.Cells.Item(26,10).Value = 1800
$ws.Cells.Item(26,11).Value = 3300
$ws.Cells.Item(26,12).Value = 3500
$ws.Cells.Item(26,13).Value = 3400
$ws.Cells.Item(26,14).Value = "`$/unidad"
$ws.Cells.Item(26,15).Value = "Región de O'Higgins"
$ws.Cells.Item(26,16).Value = 3400
$ws.Cells.Item(26,17).Value = 1
$ws.Cells.Item(26,18).Value = "Hortaliza"

# New row 27: Sandia, Segunda, date 2021-12-09 (serial 44539)
$ws.Cells.Item(27,1).Value  = 8
$ws.Cells.Item(27,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(27,3).Value  = "Coquimbo"
$ws.Cells.Item(27,4).Value  = 44539
$ws.Cells.Item(27,5).Value  = 4
$ws.Cells.Item(27,6).Value  = 100112028
$ws.Cells.Item(27,7).Value  = "Sandia"
$ws.Cells.Item(27,8).Value  = "Sin especificar"
$ws.Cells.Item(27,9).Value  = "Segunda"
$ws.Cells.Item(27,10).Value = 1000
$ws.Cells.Item(27,11).Value = 2800
$ws.Cells.Item(27,12).Value = 3000
$ws.Cells.Item(27,13).Value = 2900
$ws.Cells.Item(27,14).Value = "`$/unidad"
$ws.Cells.Item(27,15).Value = "Región de O'Higgins"
$ws.Cells.Item(27,16).Value = 2900
$ws.Cells.Item(27,17).Value = 1
$ws.Cells.Item(27,18).Value = "Hortaliza"

# Ensure the date cells keep the date-formatted style (same as other D column cells).
$ws.Range("D26").NumberFormat = $ws.Range("D28").NumberFormat
$ws.Range("D27").NumberFormat = $ws.Range("D28").NumberFormat
